# Updates cryptos list (Thu May 30 11:36:18 UTC 2024) with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.082.60"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +0.29%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.752.29"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -1.91%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.79"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.92%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.09"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -1.32%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.749.71"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -1.87%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.06%  "

# Row 9 - XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -1.87%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -3.34%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  -0.43%  "

# Row 12 - Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -2.12%  "

# Row 13 - ShibaInu
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -4.36%  "

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.66"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -1.16%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.380.08"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -2.00%  "

# Row 16 - WrappedEther
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.744.42"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -1.75%  "

# Row 17 - WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.015.80"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -0.05%  "

# Row 18 - Chainlink
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.22"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -2.61%  "

# Row 19 - Polkadot
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.07"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -5.24%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  -0.42%  "

# Row 21 - Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.82"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.43%  "

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.45"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.23%  "

# Row 23 - Polygon
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.705"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -5.04%  "

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.22"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -0.68%  "

# Row 25 - Fetch.AI
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.24"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -3.18%  "

# Row 26 - PEPE
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000137"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -9.38%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.17"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.59%  "

# Row 28 - RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.21"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -2.30%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.01%  "

# Row 30 - WrappedeETH
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.897.48"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -1.94%  "

# Row 31 - PancakeSwap
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.80"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -4.61%  "

# Row 32 - NEARProtocol
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.40"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -4.90%  "

# Row 33 - ImmutableX
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.25"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -1.71%  "

# Row 34 - EthereumClassic
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.93"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -2.79%  "

# Row 35 - Aptos
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.15"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -2.04%  "

# Row 36 - Binance-PegBSC-USD
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").NumberFormat = "General"

# Row 37 - RenzoRestakedETH
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.703.43"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -2.36%  "

# Row 38 - Hedera
$ws.Range("E38").Value = "  -4.13%  "

# Row 39 - dogwifhat
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.46"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -11.52%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  -1.09%  "

# Row 41 - Mantle
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.995"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -1.97%  "

# Row 42 - Filecoin
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.79"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -3.22%  "

# Row 43 - FirstDigitalUSD
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -0.19%  "

# Row 44 - TheGraph
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.309"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -3.38%  "

# Row 45 - USDe
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.02%  "

# Row 46 - Cosmos
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.63"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -2.00%  "

# Row 47 - Stacks
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.93"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -2.64%  "

# Row 48 - OKB
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.38"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -2.58%  "

# Row 49 - Bittensor
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "395.29"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -3.79%  "

# Row 50 - Monero
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.94"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +0.30%  "

# Row 51 - EnergySwap
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.51"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -0.03%  "

